# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (default Office palette), used by the Notes Master
#   ppt/theme/theme2.xml -> "Integral" theme, used by the (single) Slide Master / the whole deck
#
# The authored change swaps which palette is applied to the deck: the Slide
# Master's theme becomes the default "Office Theme" color scheme (and the
# Notes Master keeps/gets the "Integral" colors). Apply this the same way a
# user would from the Design tab: recolor the live theme through the
# ThemeColorScheme object exposed on a slide (it is shared by every slide
# because they all hang off the single Slide Master / theme part).

$p = $ppt.ActivePresentation

function ToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $b*65536 + $g*256 + $r
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToOleColor($officeColors[$i-1])
}
